# GMS Data Release 1
# Rename the "patient_id" field to "participant_id" in the sample field list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B8").Value = "participant_id"
